$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old sale records, keep only the header + first two rows
$ws.Range("A4:D11").EntireRow.Delete()

# Narrow the amount/payment-method columns a bit
$ws.Columns.Item(2).ColumnWidth = 6.15
$ws.Columns.Item(3).ColumnWidth = 18.15

# Row 2: refresh the sample sale record
$ws.Range("A2").Value = "2026-02-03 22:05:46"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "123.0"
$ws.Range("C2").Value = "Tarjeta de debito"

# Row 3: refresh the sample sale record
$ws.Range("A3").Value = "2026-02-05 15:15:40"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "123.0"
$ws.Range("C3").Value = "Efectivo"
